$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (columns reordered + new audit columns appended) ---
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "lang_code"
$ws.Range("E1").Value = "is_active"
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# --- Row 2: MIR / Monsieur ---
$ws.Range("A2").Value = "MIR"
$ws.Range("B2").Value = "Monsieur"
$ws.Range("C2").Value = "Titre masculin"
$ws.Range("D2").Value = "fra"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = "superadmin"
$ws.Range("G2").Value = 45079.577516400466
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = "NULL"

# --- Row 3: MRS / Madame ---
$ws.Range("A3").Value = "MRS"
$ws.Range("B3").Value = "Madame"
$ws.Range("C3").Value = "Titre fÃ©minin"
$ws.Range("D3").Value = "fra"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = "superadmin"
$ws.Range("G3").Value = 45079.577516400466
$ws.Range("G3").NumberFormat = "mm:ss.0"
$ws.Range("H3").Value = "NULL"
$ws.Range("I3").Value = "NULL"
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = "NULL"

# --- Row 4: MIS / Mademoiselle ---
$ws.Range("A4").Value = "MIS"
$ws.Range("B4").Value = "Mademoiselle"
$ws.Range("C4").Value = "Titre de femme cÃ©libataire"
$ws.Range("D4").Value = "fra"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = "superadmin"
$ws.Range("G4").Value = 45079.577516400466
$ws.Range("G4").NumberFormat = "mm:ss.0"
$ws.Range("H4").Value = "NULL"
$ws.Range("I4").Value = "NULL"
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = "NULL"

# --- Selection state ---
$ws.Range("D13").Select()
